$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = 51932000076
$ws.Range("A2").Select() | Out-Null
